$wb = $excel.ActiveWorkbook
$mismatchCount = 0

$ws = $wb.Worksheets.Item('Citywide Totals')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 2759) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L2 expected 2759 got ' + $cell.Value2) }
$cell.Value2 = 2778
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 2802) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L3 expected 2802 got ' + $cell.Value2) }
$cell.Value2 = 2822
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 749) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L4 expected 749 got ' + $cell.Value2) }
$cell.Value2 = 758
$cell = $ws.Range('L5')
if ($cell.Value2 -ne 161) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L5 expected 161 got ' + $cell.Value2) }
$cell.Value2 = 163
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 2494) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L6 expected 2494 got ' + $cell.Value2) }
$cell.Value2 = 2518
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 8965) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Citywide Totals L7 expected 8965 got ' + $cell.Value2) }
$cell.Value2 = 9039

$ws = $wb.Worksheets.Item('By Neighborhood')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 71) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L2 expected 71 got ' + $cell.Value2) }
$cell.Value2 = 73
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 301) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L7 expected 301 got ' + $cell.Value2) }
$cell.Value2 = 304
$cell = $ws.Range('L8')
if ($cell.Value2 -ne 566) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L8 expected 566 got ' + $cell.Value2) }
$cell.Value2 = 573
$cell = $ws.Range('L10')
if ($cell.Value2 -ne 57) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L10 expected 57 got ' + $cell.Value2) }
$cell.Value2 = 58
$cell = $ws.Range('L11')
if ($cell.Value2 -ne 158) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L11 expected 158 got ' + $cell.Value2) }
$cell.Value2 = 160
$cell = $ws.Range('L14')
if ($cell.Value2 -ne 41) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L14 expected 41 got ' + $cell.Value2) }
$cell.Value2 = 43
$cell = $ws.Range('L18')
if ($cell.Value2 -ne 63) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L18 expected 63 got ' + $cell.Value2) }
$cell.Value2 = 64
$cell = $ws.Range('L19')
if ($cell.Value2 -ne 249) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L19 expected 249 got ' + $cell.Value2) }
$cell.Value2 = 252
$cell = $ws.Range('L20')
if ($cell.Value2 -ne 230) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L20 expected 230 got ' + $cell.Value2) }
$cell.Value2 = 231
$cell = $ws.Range('L22')
if ($cell.Value2 -ne 28) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L22 expected 28 got ' + $cell.Value2) }
$cell.Value2 = 29
$cell = $ws.Range('L24')
if ($cell.Value2 -ne 20) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L24 expected 20 got ' + $cell.Value2) }
$cell.Value2 = 21
$cell = $ws.Range('L26')
if ($cell.Value2 -ne 8) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L26 expected 8 got ' + $cell.Value2) }
$cell.Value2 = 9
$cell = $ws.Range('L29')
if ($cell.Value2 -ne 481) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L29 expected 481 got ' + $cell.Value2) }
$cell.Value2 = 485
$cell = $ws.Range('L30')
if ($cell.Value2 -ne 42) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L30 expected 42 got ' + $cell.Value2) }
$cell.Value2 = 44
$cell = $ws.Range('L33')
if ($cell.Value2 -ne 409) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L33 expected 409 got ' + $cell.Value2) }
$cell.Value2 = 411
$cell = $ws.Range('L37')
if ($cell.Value2 -ne 329) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L37 expected 329 got ' + $cell.Value2) }
$cell.Value2 = 332
$cell = $ws.Range('L42')
if ($cell.Value2 -ne 298) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L42 expected 298 got ' + $cell.Value2) }
$cell.Value2 = 301
$cell = $ws.Range('L44')
if ($cell.Value2 -ne 67) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L44 expected 67 got ' + $cell.Value2) }
$cell.Value2 = 68
$cell = $ws.Range('L47')
if ($cell.Value2 -ne 68) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L47 expected 68 got ' + $cell.Value2) }
$cell.Value2 = 69
$cell = $ws.Range('L48')
if ($cell.Value2 -ne 119) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L48 expected 119 got ' + $cell.Value2) }
$cell.Value2 = 120
$cell = $ws.Range('L52')
if ($cell.Value2 -ne 180) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L52 expected 180 got ' + $cell.Value2) }
$cell.Value2 = 182
$cell = $ws.Range('L55')
if ($cell.Value2 -ne 86) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L55 expected 86 got ' + $cell.Value2) }
$cell.Value2 = 87
$cell = $ws.Range('L60')
if ($cell.Value2 -ne 55) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L60 expected 55 got ' + $cell.Value2) }
$cell.Value2 = 56
$cell = $ws.Range('L63')
if ($cell.Value2 -ne 30) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L63 expected 30 got ' + $cell.Value2) }
$cell.Value2 = 29
$cell = $ws.Range('L64')
if ($cell.Value2 -ne 55) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L64 expected 55 got ' + $cell.Value2) }
$cell.Value2 = 56
$cell = $ws.Range('L65')
if ($cell.Value2 -ne 163) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L65 expected 163 got ' + $cell.Value2) }
$cell.Value2 = 165
$cell = $ws.Range('L67')
if ($cell.Value2 -ne 331) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L67 expected 331 got ' + $cell.Value2) }
$cell.Value2 = 335
$cell = $ws.Range('L72')
if ($cell.Value2 -ne 39) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L72 expected 39 got ' + $cell.Value2) }
$cell.Value2 = 41
$cell = $ws.Range('L75')
if ($cell.Value2 -ne 34) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L75 expected 34 got ' + $cell.Value2) }
$cell.Value2 = 35
$cell = $ws.Range('L76')
if ($cell.Value2 -ne 112) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L76 expected 112 got ' + $cell.Value2) }
$cell.Value2 = 115
$cell = $ws.Range('L79')
if ($cell.Value2 -ne 240) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L79 expected 240 got ' + $cell.Value2) }
$cell.Value2 = 241
$cell = $ws.Range('L83')
if ($cell.Value2 -ne 211) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L83 expected 211 got ' + $cell.Value2) }
$cell.Value2 = 213
$cell = $ws.Range('L84')
if ($cell.Value2 -ne 93) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L84 expected 93 got ' + $cell.Value2) }
$cell.Value2 = 94
$cell = $ws.Range('L85')
if ($cell.Value2 -ne 459) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L85 expected 459 got ' + $cell.Value2) }
$cell.Value2 = 461
$cell = $ws.Range('L86')
if ($cell.Value2 -ne 64) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L86 expected 64 got ' + $cell.Value2) }
$cell.Value2 = 65
$cell = $ws.Range('L89')
if ($cell.Value2 -ne 114) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L89 expected 114 got ' + $cell.Value2) }
$cell.Value2 = 115
$cell = $ws.Range('L91')
if ($cell.Value2 -ne 127) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L91 expected 127 got ' + $cell.Value2) }
$cell.Value2 = 128
$cell = $ws.Range('L94')
if ($cell.Value2 -ne 109) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L94 expected 109 got ' + $cell.Value2) }
$cell.Value2 = 111
$cell = $ws.Range('L96')
if ($cell.Value2 -ne 88) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L96 expected 88 got ' + $cell.Value2) }
$cell.Value2 = 90
$cell = $ws.Range('L97')
if ($cell.Value2 -ne 80) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L97 expected 80 got ' + $cell.Value2) }
$cell.Value2 = 81
$cell = $ws.Range('L99')
if ($cell.Value2 -ne 149) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L99 expected 149 got ' + $cell.Value2) }
$cell.Value2 = 150
$cell = $ws.Range('L100')
if ($cell.Value2 -ne 14) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L100 expected 14 got ' + $cell.Value2) }
$cell.Value2 = 15
$cell = $ws.Range('L101')
if ($cell.Value2 -ne 8965) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH By Neighborhood L101 expected 8965 got ' + $cell.Value2) }
$cell.Value2 = 9039

$ws = $wb.Worksheets.Item('Austin')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 190) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Austin L3 expected 190 got ' + $cell.Value2) }
$cell.Value2 = 193
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 40) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Austin L4 expected 40 got ' + $cell.Value2) }
$cell.Value2 = 41
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 155) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Austin L6 expected 155 got ' + $cell.Value2) }
$cell.Value2 = 158
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 566) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Austin L7 expected 566 got ' + $cell.Value2) }
$cell.Value2 = 573

$ws = $wb.Worksheets.Item('South Chicago')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 88) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Chicago L3 expected 88 got ' + $cell.Value2) }
$cell.Value2 = 89
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 45) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Chicago L6 expected 45 got ' + $cell.Value2) }
$cell.Value2 = 46
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 211) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Chicago L7 expected 211 got ' + $cell.Value2) }
$cell.Value2 = 213

$ws = $wb.Worksheets.Item('Garfield Park')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 113) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Garfield Park L2 expected 113 got ' + $cell.Value2) }
$cell.Value2 = 114
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 140) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Garfield Park L6 expected 140 got ' + $cell.Value2) }
$cell.Value2 = 141
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 409) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Garfield Park L7 expected 409 got ' + $cell.Value2) }
$cell.Value2 = 411

$ws = $wb.Worksheets.Item('Grand Crossing')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 95) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Grand Crossing L3 expected 95 got ' + $cell.Value2) }
$cell.Value2 = 96
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 21) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Grand Crossing L4 expected 21 got ' + $cell.Value2) }
$cell.Value2 = 22
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 105) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Grand Crossing L6 expected 105 got ' + $cell.Value2) }
$cell.Value2 = 106
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 329) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Grand Crossing L7 expected 329 got ' + $cell.Value2) }
$cell.Value2 = 332

$ws = $wb.Worksheets.Item('New City')
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 7) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH New City L4 expected 7 got ' + $cell.Value2) }
$cell.Value2 = 8
$cell = $ws.Range('L5')
if ($cell.Value2 -ne 2) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH New City L5 expected 2 got ' + $cell.Value2) }
$cell.Value2 = 3
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 163) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH New City L7 expected 163 got ' + $cell.Value2) }
$cell.Value2 = 165

$ws = $wb.Worksheets.Item('Woodlawn')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 36) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Woodlawn L6 expected 36 got ' + $cell.Value2) }
$cell.Value2 = 37
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 149) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Woodlawn L7 expected 149 got ' + $cell.Value2) }
$cell.Value2 = 150

$ws = $wb.Worksheets.Item('Fuller Park')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 19) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Fuller Park L6 expected 19 got ' + $cell.Value2) }
$cell.Value2 = 21
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 42) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Fuller Park L7 expected 42 got ' + $cell.Value2) }
$cell.Value2 = 44

$ws = $wb.Worksheets.Item('North Lawndale')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 97) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH North Lawndale L2 expected 97 got ' + $cell.Value2) }
$cell.Value2 = 99
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 120) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH North Lawndale L3 expected 120 got ' + $cell.Value2) }
$cell.Value2 = 122
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 331) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH North Lawndale L7 expected 331 got ' + $cell.Value2) }
$cell.Value2 = 335

$ws = $wb.Worksheets.Item('South Deering')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 20) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Deering L6 expected 20 got ' + $cell.Value2) }
$cell.Value2 = 21
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 93) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Deering L7 expected 93 got ' + $cell.Value2) }
$cell.Value2 = 94

$ws = $wb.Worksheets.Item('Englewood')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 154) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Englewood L2 expected 154 got ' + $cell.Value2) }
$cell.Value2 = 156
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 19) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Englewood L4 expected 19 got ' + $cell.Value2) }
$cell.Value2 = 20
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 123) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Englewood L6 expected 123 got ' + $cell.Value2) }
$cell.Value2 = 124
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 481) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Englewood L7 expected 481 got ' + $cell.Value2) }
$cell.Value2 = 485

$ws = $wb.Worksheets.Item('Lake View')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 14) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Lake View L2 expected 14 got ' + $cell.Value2) }
$cell.Value2 = 15
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 119) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Lake View L7 expected 119 got ' + $cell.Value2) }
$cell.Value2 = 120

$ws = $wb.Worksheets.Item('Chatham')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 78) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chatham L3 expected 78 got ' + $cell.Value2) }
$cell.Value2 = 79
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 9) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chatham L4 expected 9 got ' + $cell.Value2) }
$cell.Value2 = 10
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 76) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chatham L6 expected 76 got ' + $cell.Value2) }
$cell.Value2 = 77
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 249) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chatham L7 expected 249 got ' + $cell.Value2) }
$cell.Value2 = 252

$ws = $wb.Worksheets.Item('Irving Park')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 29) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Irving Park L2 expected 29 got ' + $cell.Value2) }
$cell.Value2 = 30
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 67) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Irving Park L7 expected 67 got ' + $cell.Value2) }
$cell.Value2 = 68

$ws = $wb.Worksheets.Item('River North')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 22) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH River North L2 expected 22 got ' + $cell.Value2) }
$cell.Value2 = 23
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 54) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH River North L6 expected 54 got ' + $cell.Value2) }
$cell.Value2 = 56
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 112) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH River North L7 expected 112 got ' + $cell.Value2) }
$cell.Value2 = 115

$ws = $wb.Worksheets.Item('Bridgeport')
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 4) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Bridgeport L4 expected 4 got ' + $cell.Value2) }
$cell.Value2 = 6
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 41) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Bridgeport L7 expected 41 got ' + $cell.Value2) }
$cell.Value2 = 43

$ws = $wb.Worksheets.Item('Humboldt Park')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 85) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Humboldt Park L2 expected 85 got ' + $cell.Value2) }
$cell.Value2 = 87
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 90) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Humboldt Park L3 expected 90 got ' + $cell.Value2) }
$cell.Value2 = 91
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 298) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Humboldt Park L7 expected 298 got ' + $cell.Value2) }
$cell.Value2 = 301

$ws = $wb.Worksheets.Item('Avondale')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 14) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Avondale L3 expected 14 got ' + $cell.Value2) }
$cell.Value2 = 15
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 57) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Avondale L7 expected 57 got ' + $cell.Value2) }
$cell.Value2 = 58

$ws = $wb.Worksheets.Item('Lower West Side')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 30) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Lower West Side L3 expected 30 got ' + $cell.Value2) }
$cell.Value2 = 31
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 86) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Lower West Side L7 expected 86 got ' + $cell.Value2) }
$cell.Value2 = 87

$ws = $wb.Worksheets.Item('Dunning')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 9) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Dunning L2 expected 9 got ' + $cell.Value2) }
$cell.Value2 = 10
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 20) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Dunning L7 expected 20 got ' + $cell.Value2) }
$cell.Value2 = 21

$ws = $wb.Worksheets.Item('West Ridge')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 22) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Ridge L3 expected 22 got ' + $cell.Value2) }
$cell.Value2 = 23
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 20) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Ridge L6 expected 20 got ' + $cell.Value2) }
$cell.Value2 = 21
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 88) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Ridge L7 expected 88 got ' + $cell.Value2) }
$cell.Value2 = 90

$ws = $wb.Worksheets.Item('Washington Park')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 49) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Washington Park L2 expected 49 got ' + $cell.Value2) }
$cell.Value2 = 50
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 127) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Washington Park L7 expected 127 got ' + $cell.Value2) }
$cell.Value2 = 128

$ws = $wb.Worksheets.Item('Roseland')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 77) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Roseland L2 expected 77 got ' + $cell.Value2) }
$cell.Value2 = 78
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 240) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Roseland L7 expected 240 got ' + $cell.Value2) }
$cell.Value2 = 241

$ws = $wb.Worksheets.Item('Near South Side')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 19) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Near South Side L2 expected 19 got ' + $cell.Value2) }
$cell.Value2 = 20
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 55) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Near South Side L7 expected 55 got ' + $cell.Value2) }
$cell.Value2 = 56

$ws = $wb.Worksheets.Item('Chicago Lawn')
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 20) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chicago Lawn L4 expected 20 got ' + $cell.Value2) }
$cell.Value2 = 21
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 230) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Chicago Lawn L7 expected 230 got ' + $cell.Value2) }
$cell.Value2 = 231

$ws = $wb.Worksheets.Item('Calumet Heights')
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 6) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Calumet Heights L4 expected 6 got ' + $cell.Value2) }
$cell.Value2 = 7
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 63) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Calumet Heights L7 expected 63 got ' + $cell.Value2) }
$cell.Value2 = 64

$ws = $wb.Worksheets.Item('Wrigleyville')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 5) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Wrigleyville L3 expected 5 got ' + $cell.Value2) }
$cell.Value2 = 6
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 14) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Wrigleyville L7 expected 14 got ' + $cell.Value2) }
$cell.Value2 = 15

$ws = $wb.Worksheets.Item('Auburn Gresham')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 92) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Auburn Gresham L2 expected 92 got ' + $cell.Value2) }
$cell.Value2 = 93
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 95) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Auburn Gresham L3 expected 95 got ' + $cell.Value2) }
$cell.Value2 = 96
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 82) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Auburn Gresham L6 expected 82 got ' + $cell.Value2) }
$cell.Value2 = 83
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 301) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Auburn Gresham L7 expected 301 got ' + $cell.Value2) }
$cell.Value2 = 304

$ws = $wb.Worksheets.Item('West Loop')
$cell = $ws.Range('L4')
if ($cell.Value2 -ne 15) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Loop L4 expected 15 got ' + $cell.Value2) }
$cell.Value2 = 16
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 35) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Loop L6 expected 35 got ' + $cell.Value2) }
$cell.Value2 = 36
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 109) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Loop L7 expected 109 got ' + $cell.Value2) }
$cell.Value2 = 111

$ws = $wb.Worksheets.Item('Kenwood')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 23) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Kenwood L3 expected 23 got ' + $cell.Value2) }
$cell.Value2 = 24
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 68) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Kenwood L7 expected 68 got ' + $cell.Value2) }
$cell.Value2 = 69

$ws = $wb.Worksheets.Item('East Village')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 5) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH East Village L6 expected 5 got ' + $cell.Value2) }
$cell.Value2 = 6
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 8) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH East Village L7 expected 8 got ' + $cell.Value2) }
$cell.Value2 = 9

$ws = $wb.Worksheets.Item('Belmont Cragin')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 50) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Belmont Cragin L3 expected 50 got ' + $cell.Value2) }
$cell.Value2 = 51
$cell = $ws.Range('L5')
if ($cell.Value2 -ne 1) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Belmont Cragin L5 expected 1 got ' + $cell.Value2) }
$cell.Value2 = 2
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 158) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Belmont Cragin L7 expected 158 got ' + $cell.Value2) }
$cell.Value2 = 160

$ws = $wb.Worksheets.Item('Albany Park')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 21) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Albany Park L6 expected 21 got ' + $cell.Value2) }
$cell.Value2 = 23
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 71) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Albany Park L7 expected 71 got ' + $cell.Value2) }
$cell.Value2 = 73

$ws = $wb.Worksheets.Item('West Town')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 42) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Town L6 expected 42 got ' + $cell.Value2) }
$cell.Value2 = 43
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 80) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH West Town L7 expected 80 got ' + $cell.Value2) }
$cell.Value2 = 81

$ws = $wb.Worksheets.Item('Uptown')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 30) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Uptown L3 expected 30 got ' + $cell.Value2) }
$cell.Value2 = 31
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 114) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Uptown L7 expected 114 got ' + $cell.Value2) }
$cell.Value2 = 115

$ws = $wb.Worksheets.Item('Streeterville')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 11) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Streeterville L3 expected 11 got ' + $cell.Value2) }
$cell.Value2 = 12
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 64) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Streeterville L7 expected 64 got ' + $cell.Value2) }
$cell.Value2 = 65

$ws = $wb.Worksheets.Item('Pullman')
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 2) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Pullman L6 expected 2 got ' + $cell.Value2) }
$cell.Value2 = 3
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 34) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Pullman L7 expected 34 got ' + $cell.Value2) }
$cell.Value2 = 35

$ws = $wb.Worksheets.Item('Morgan Park')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 21) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Morgan Park L3 expected 21 got ' + $cell.Value2) }
$cell.Value2 = 22
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 55) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Morgan Park L7 expected 55 got ' + $cell.Value2) }
$cell.Value2 = 56

$ws = $wb.Worksheets.Item('South Shore')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 133) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Shore L2 expected 133 got ' + $cell.Value2) }
$cell.Value2 = 134
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 93) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Shore L6 expected 93 got ' + $cell.Value2) }
$cell.Value2 = 94
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 459) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH South Shore L7 expected 459 got ' + $cell.Value2) }
$cell.Value2 = 461

$ws = $wb.Worksheets.Item('Clearing')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 10) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Clearing L2 expected 10 got ' + $cell.Value2) }
$cell.Value2 = 11
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 28) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Clearing L7 expected 28 got ' + $cell.Value2) }
$cell.Value2 = 29

$ws = $wb.Worksheets.Item('Old Town')
$cell = $ws.Range('L3')
if ($cell.Value2 -ne 7) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Old Town L3 expected 7 got ' + $cell.Value2) }
$cell.Value2 = 8
$cell = $ws.Range('L6')
if ($cell.Value2 -ne 11) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Old Town L6 expected 11 got ' + $cell.Value2) }
$cell.Value2 = 12
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 39) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Old Town L7 expected 39 got ' + $cell.Value2) }
$cell.Value2 = 41

$ws = $wb.Worksheets.Item('Little Village')
$cell = $ws.Range('L2')
if ($cell.Value2 -ne 61) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Little Village L2 expected 61 got ' + $cell.Value2) }
$cell.Value2 = 63
$cell = $ws.Range('L7')
if ($cell.Value2 -ne 180) { $mismatchCount = $mismatchCount + 1; Write-Host ('MISMATCH Little Village L7 expected 180 got ' + $cell.Value2) }
$cell.Value2 = 182

Write-Host ("Total mismatches: " + $mismatchCount)
Write-Host "Done applying 2025-06-09 data update."